# Update column F (dSF) values for specific rows to match the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F2"  = -5
    "F4"  = 3
    "F6"  = -1
    "F7"  = 3
    "F9"  = 1
    "F10" = 1
    "F13" = -3
    "F14" = -9
    "F29" = -3
    "F30" = -5
    "F31" = 0
    "F32" = 0
    "F36" = 2
    "F37" = 3
    "F38" = 2
    "F42" = 0
    "F43" = 2
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
